# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates (currentAveragePrice*,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ -- columns H..N)
# to the specific Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# that changed in this refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 1055.8182
$ws.Range("I12").Value = 1148.1
$ws.Range("J12").Value = 133
$ws.Range("K12").Value = 1148.1
$ws.Range("L12").Value = 133
$ws.Range("M12").Value = -978.0999999999999
$ws.Range("N12").Value = -473

# Row 19
$ws.Range("H19").Value = 1948.4546
$ws.Range("I19").Value = 1107.8
$ws.Range("J19").Value = 2649
$ws.Range("K19").Value = 1107.8
$ws.Range("L19").Value = 2649
$ws.Range("M19").Value = -932.8
$ws.Range("N19").Value = -2999

# Row 40
$ws.Range("H40").Value = 5119.643
$ws.Range("I40").Value = 2425
$ws.Range("J40").Value = 5326.923
$ws.Range("K40").Value = 2425
$ws.Range("L40").Value = 5326.923
$ws.Range("M40").Value = -2250
$ws.Range("N40").Value = -5676.923

# Row 64
$ws.Range("H64").Value = 3898.9285
$ws.Range("I64").Value = 3749
$ws.Range("K64").Value = 3749
$ws.Range("M64").Value = -3501

# Row 67
$ws.Range("H67").Value = 3898.9285
$ws.Range("I67").Value = 3749
$ws.Range("K67").Value = 3749
$ws.Range("M67").Value = -2891

# Row 76
$ws.Range("H76").Value = 4890.0557
$ws.Range("I76").Value = 3000.3845
$ws.Range("J76").Value = 9803.200000000001
$ws.Range("K76").Value = 3000.3845
$ws.Range("L76").Value = 9803.200000000001
$ws.Range("M76").Value = -2685.3845
$ws.Range("N76").Value = -10433.2

# Row 79
$ws.Range("H79").Value = 4890.0557
$ws.Range("I79").Value = 3000.3845
$ws.Range("J79").Value = 9803.200000000001
$ws.Range("K79").Value = 3000.3845
$ws.Range("L79").Value = 9803.200000000001
$ws.Range("M79").Value = -1908.3845
$ws.Range("N79").Value = -11987.2

# Row 106
$ws.Range("H106").Value = 8515.916999999999
$ws.Range("I106").Value = 4998.875
$ws.Range("J106").Value = 15550
$ws.Range("K106").Value = 4998.875
$ws.Range("L106").Value = 15550
$ws.Range("M106").Value = -4367.875
$ws.Range("N106").Value = -16812

# Row 107
$ws.Range("H107").Value = 3376.1482
$ws.Range("I107").Value = 2582.923
$ws.Range("K107").Value = 2582.923
$ws.Range("M107").Value = -662.9229999999998

# Row 136
$ws.Range("H136").Value = 90000
$ws.Range("J136").Value = 90000
$ws.Range("L136").Value = 90000
$ws.Range("N136").Value = -100200

# Row 137
$ws.Range("H137").Value = 5498.5356
$ws.Range("I137").Value = 2009.1305
$ws.Range("J137").Value = 21549.8
$ws.Range("K137").Value = 6027.3915
$ws.Range("L137").Value = 64649.39999999999
$ws.Range("M137").Value = -3477.3915
$ws.Range("N137").Value = -69749.39999999999

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 201682.4
$ws.Range("I32").Value = 221252.1
$ws.Range("J32").Value = 21641.4
$ws.Range("K32").Value = 221252.1
$ws.Range("L32").Value = 21641.4
$ws.Range("M32").Value = -220965.1
$ws.Range("N32").Value = -22215.4

# Row 63
$ws.Range("H63").Value = 4618.625
$ws.Range("I63").Value = 1949.5
$ws.Range("K63").Value = 1949.5
$ws.Range("M63").Value = -1263.5

# Row 66
$ws.Range("H66").Value = 4618.625
$ws.Range("I66").Value = 1949.5
$ws.Range("K66").Value = 9747.5
$ws.Range("M66").Value = -6315.5

# Row 74
$ws.Range("H74").Value = 9084.441000000001
$ws.Range("I74").Value = 6153.76
$ws.Range("K74").Value = 6153.76
$ws.Range("M74").Value = -5279.76

# Row 77
$ws.Range("H77").Value = 9084.441000000001
$ws.Range("I77").Value = 6153.76
$ws.Range("K77").Value = 30768.8
$ws.Range("M77").Value = -26400.8

# Row 88
$ws.Range("H88").Value = 2290.889
$ws.Range("J88").Value = 2290.889
$ws.Range("L88").Value = 2290.889
$ws.Range("N88").Value = -3102.889

# Row 91
$ws.Range("H91").Value = 2290.889
$ws.Range("J91").Value = 2290.889
$ws.Range("L91").Value = 2290.889
$ws.Range("N91").Value = -5098.889

# Row 110
$ws.Range("H110").Value = 2546.4546
$ws.Range("I110").Value = 8011
$ws.Range("K110").Value = 8011
$ws.Range("M110").Value = -5966

# Row 132
$ws.Range("H132").Value = 6833.7427
$ws.Range("I132").Value = 5743.375
$ws.Range("J132").Value = 7751.9473
$ws.Range("K132").Value = 17230.125
$ws.Range("L132").Value = 23255.8419
$ws.Range("M132").Value = -14700.125
$ws.Range("N132").Value = -28315.8419

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 3377.0344
$ws.Range("I86").Value = 3538.4
$ws.Range("J86").Value = 3018.4443
$ws.Range("K86").Value = 3538.4
$ws.Range("L86").Value = 3018.4443
$ws.Range("M86").Value = -2415.4
$ws.Range("N86").Value = -5264.4443

# Row 89
$ws.Range("H89").Value = 3377.0344
$ws.Range("I89").Value = 3538.4
$ws.Range("J89").Value = 3018.4443
$ws.Range("K89").Value = 17692
$ws.Range("L89").Value = 15092.2215
$ws.Range("M89").Value = -12076
$ws.Range("N89").Value = -26324.2215

# Row 105
$ws.Range("H105").Value = 10486.5
$ws.Range("I105").Value = 1004.5
$ws.Range("K105").Value = 1004.5
$ws.Range("M105").Value = 742.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 5750
$ws.Range("I31").Value = 5750
$ws.Range("K31").Value = 5750
$ws.Range("M31").Value = -5455

# Row 34
$ws.Range("H34").Value = 5750
$ws.Range("I34").Value = 5750
$ws.Range("K34").Value = 5750
$ws.Range("M34").Value = -5548

# Row 62
$ws.Range("H62").Value = 3962.889
$ws.Range("I62").Value = 3566.3333
$ws.Range("K62").Value = 3566.3333
$ws.Range("M62").Value = -2942.3333

# Row 65
$ws.Range("H65").Value = 3962.889
$ws.Range("I65").Value = 3566.3333
$ws.Range("K65").Value = 17831.6665
$ws.Range("M65").Value = -14711.6665

# Row 105
$ws.Range("H105").Value = 12315.077
$ws.Range("I105").Value = 13091.333
$ws.Range("K105").Value = 13091.333
$ws.Range("M105").Value = -11344.333

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 2001382.5
$ws.Range("J4").Value = 2662.6667
$ws.Range("L4").Value = 7988.000100000001
$ws.Range("N4").Value = -8212.000100000001

# Row 37
$ws.Range("H37").Value = 50000
$ws.Range("J37").Value = 50000
$ws.Range("L37").Value = 150000
$ws.Range("N37").Value = -150224

# Row 128
$ws.Range("H128").Value = 320294.88
$ws.Range("I128").Value = 320294.88
$ws.Range("K128").Value = 960884.64
$ws.Range("M128").Value = -955904.64

# Row 131
$ws.Range("H131").Value = 10799.714
$ws.Range("J131").Value = 10972.454
$ws.Range("L131").Value = 32917.362
$ws.Range("N131").Value = -42997.362

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 59
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

# Row 70
$ws.Range("H70").Value = 19103.588
$ws.Range("I70").Value = 22661.916
$ws.Range("J70").Value = 10563.6
$ws.Range("K70").Value = 22661.916
$ws.Range("L70").Value = 10563.6
$ws.Range("M70").Value = -22391.916
$ws.Range("N70").Value = -11103.6

# Row 73
$ws.Range("H73").Value = 19103.588
$ws.Range("I73").Value = 22661.916
$ws.Range("J73").Value = 10563.6
$ws.Range("K73").Value = 22661.916
$ws.Range("L73").Value = 10563.6
$ws.Range("M73").Value = -21725.916
$ws.Range("N73").Value = -12435.6

# Row 80
$ws.Range("H80").Value = 3138
$ws.Range("I80").Value = 2656.2222
$ws.Range("J80").Value = 4583.3335
$ws.Range("K80").Value = 2656.2222
$ws.Range("L80").Value = 4583.3335
$ws.Range("M80").Value = -1658.2222
$ws.Range("N80").Value = -6579.3335

# Row 83
$ws.Range("H83").Value = 3138
$ws.Range("I83").Value = 2656.2222
$ws.Range("J83").Value = 4583.3335
$ws.Range("K83").Value = 13281.111
$ws.Range("L83").Value = 22916.6675
$ws.Range("M83").Value = -8289.111000000001
$ws.Range("N83").Value = -32900.6675

# Row 126
$ws.Range("H126").Value = 2487
$ws.Range("I126").Value = 2328.5557
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 6985.6671
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -4515.6671
$ws.Range("N126").Value = -14540

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 50000
$ws.Range("I40").Value = 50000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 50000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -49864
$ws.Range("N40").ClearContents()

# Row 108
$ws.Range("H108").Value = 41546
$ws.Range("J108").Value = 41546
$ws.Range("L108").Value = 41546
$ws.Range("N108").Value = -49226

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 143977.62
$ws.Range("I122").Value = 2233.3333
$ws.Range("J122").Value = 229024.2
$ws.Range("K122").Value = 6699.999899999999
$ws.Range("L122").Value = 687072.6000000001
$ws.Range("M122").Value = -4249.999899999999
$ws.Range("N122").Value = -691972.6000000001

# Row 132
$ws.Range("H132").Value = 2902.3076
$ws.Range("I132").Value = 2466
$ws.Range("K132").Value = 7398
$ws.Range("M132").Value = -4868
